$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (offsets chosen so the stored OOXML "width"
# attribute lands exactly on the target value after Excel's internal
# pixel-rounding of ColumnWidth)
$ws.Columns.Item(3).ColumnWidth = 65.165
$ws.Columns.Item(4).ColumnWidth = 59.165
$ws.Columns.Item(8).ColumnWidth = 56.165

# Column A holds opportunity IDs as TEXT (numeric-looking strings). Force
# text formatting first so the new numeric-looking IDs below don't get
# auto-converted to Number cells.
$ws.Range("A2:A14").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "1328351"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328351"
$ws.Range("C2").Value = "Test Engineer"
$ws.Range("D2").Value = "Hamburg, Deutschland"
$ws.Range("F2").Value = "15 applicants"
$ws.Range("H2").Value = "akeno"

# Row 3
$ws.Range("A3").Value = "1328344"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328344"
$ws.Range("C3").Value = "Business Development Outbound Specialist"
$ws.Range("D3").Value = "Београд, Србија"
$ws.Range("F3").Value = "5 applicants"
$ws.Range("H3").Value = "Native Teams"

# Row 4
$ws.Range("A4").Value = "1328339"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328339"
$ws.Range("C4").Value = "IT Research Interns (Duplicated)"
$ws.Range("D4").Value = "Aronj, Uttar Pradesh, India"
$ws.Range("E4").Value = "No"
# (E4's yellow "premium" highlight (style index 3) is removed below, after
# it has been used as the format source for E9's new highlight.)
$ws.Range("F4").Value = "3 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "FS University"

# Row 5
$ws.Range("A5").Value = "1328306"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328306"
$ws.Range("C5").Value = "Automotive Maintenance Technician"
$ws.Range("D5").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("H5").Value = "Automotive fleet services"

# Row 6
$ws.Range("A6").Value = "1328250"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1328250"
$ws.Range("C6").Value = "Digital Marketing Intern"
$ws.Range("D6").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F6").Value = "2 applicants"
$ws.Range("H6").Value = "Requisite Technologies Pvt Ltd"

# Row 7
$ws.Range("A7").Value = "1328179"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1328179"
$ws.Range("C7").Value = "Social Media Executive"
$ws.Range("D7").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F7").Value = "1 applicant"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "Leadmagnet private limited"

# Row 8
$ws.Range("A8").Value = "1326536"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1326536"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("F8").Value = "12 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Egypt holiday travel"

# Row 9
$ws.Range("A9").Value = "1326481"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1326481"
$ws.Range("C9").Value = "Global Duty Billing Data Analytics Expert"
$ws.Range("D9").Value = "Maastricht, Netherlands"
# Copy the "premium highlight" format (yellow fill, style index 3) from
# E4 onto E9, then set E9's value and finally strip E4's own highlight.
$ws.Range("E4").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "Yes"
$ws.Range("E4").ClearFormats()
$ws.Range("F9").Value = "206 applicants"
$ws.Range("H9").Value = "DHL Group"

# Row 10
$ws.Range("A10").Value = "1323019"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1323019"
$ws.Range("C10").Value = "Project Specialist"
$ws.Range("D10").Value = "İstanbul, Türkiye"
$ws.Range("F10").Value = "156 applicants"
$ws.Range("H10").Value = "Nabulu"

# Row 11
$ws.Range("A11").Value = "1320868"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1320868"
$ws.Range("C11").Value = "Accelerate Romania|Data Labeling Specialist (SERBIAN Speackers)"
$ws.Range("D11").Value = "Bucharest, Romania"
$ws.Range("E11").Value = "No"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = "7 applicants"
$ws.Range("G11").Value = "9 - 12 Weeks"
$ws.Range("H11").Value = "RepsMate"

# Row 12
$ws.Range("A12").Value = "1317292"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1317292"
$ws.Range("C12").Value = "[Impact Florianópolis]- Social Media"
$ws.Range("D12").Value = "São Miguel do Oeste, SC, 89900-000, Brasil"
$ws.Range("F12").Value = "80 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "KNN Idiomas"

# Row 13
$ws.Range("A13").Value = "1315265"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1315265"
$ws.Range("C13").Value = "Sales Attendant"
$ws.Range("D13").Value = "Denizli, Kumkısık, Denizli, Türkiye"
$ws.Range("F13").Value = "66 applicants"
$ws.Range("H13").Value = "COTTON CASTLE TEKSTİL SANAYİ VE TİCARET ANONİM ŞİRKETİ"

# Row 14
$ws.Range("A14").Value = "1289379"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1289379"
$ws.Range("C14").Value = "Medical Advisor Portuguese Speaker"
$ws.Range("D14").Value = "İstanbul, Türkiye"
$ws.Range("H14").Value = "International Plus"
